$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 3175.2666  # H18
$ws.Cells.Item(18, 9).Value = 3316.2727  # I18
$ws.Cells.Item(18, 10).Value = 2787.5  # J18
$ws.Cells.Item(18, 11).Value = 3316.2727  # K18
$ws.Cells.Item(18, 12).Value = 2787.5  # L18
$ws.Cells.Item(18, 13).Value = -3032.2727  # M18
$ws.Cells.Item(18, 14).Value = -3355.5  # N18
$ws.Cells.Item(41, 8).Value = 855.9091  # H41
$ws.Cells.Item(41, 10).Value = 899.5833  # J41
$ws.Cells.Item(41, 12).Value = 899.5833  # L41
$ws.Cells.Item(41, 14).Value = -1779.5833  # N41
$ws.Cells.Item(69, 8).Value = 45484216  # H69
$ws.Cells.Item(69, 9).Value = 170341  # I69
$ws.Cells.Item(69, 11).Value = 511023  # K69
$ws.Cells.Item(69, 13).Value = -510149  # M69
$ws.Cells.Item(72, 8).Value = 45484216  # H72
$ws.Cells.Item(72, 9).Value = 170341  # I72
$ws.Cells.Item(72, 11).Value = 1533069  # K72
$ws.Cells.Item(72, 13).Value = -1528701  # M72
$ws.Cells.Item(92, 8).Value = 861.6  # H92
$ws.Cells.Item(92, 9).Value = 882.5  # I92
$ws.Cells.Item(92, 10).Value = 778  # J92
$ws.Cells.Item(92, 11).Value = 882.5  # K92
$ws.Cells.Item(92, 12).Value = 778  # L92
$ws.Cells.Item(92, 13).Value = 365.5  # M92
$ws.Cells.Item(92, 14).Value = -3274  # N92
$ws.Cells.Item(93, 8).Value = 34500  # H93
$ws.Cells.Item(93, 10).Value = 34500  # J93
$ws.Cells.Item(93, 12).Value = 34500  # L93
$ws.Cells.Item(93, 14).Value = -39492  # N93

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 78.8  # H4
$ws.Cells.Item(4, 9).Value = 78.8  # I4
$ws.Cells.Item(4, 10).Value = 0  # J4
$ws.Cells.Item(4, 11).Value = 78.8  # K4
$ws.Cells.Item(4, 12).Value = 0  # L4
$ws.Cells.Item(4, 13).Value = 37.2  # M4
$ws.Cells.Item(4, 14).ClearContents()  # N4
$ws.Cells.Item(32, 8).Value = 4413.234  # H32
$ws.Cells.Item(32, 9).Value = 2647.9106  # I32
$ws.Cells.Item(32, 10).Value = 9120.762000000001  # J32
$ws.Cells.Item(32, 11).Value = 2647.9106  # K32
$ws.Cells.Item(32, 12).Value = 9120.762000000001  # L32
$ws.Cells.Item(32, 13).Value = -2360.9106  # M32
$ws.Cells.Item(32, 14).Value = -9694.762000000001  # N32
$ws.Cells.Item(45, 8).Value = 7195936  # H45
$ws.Cells.Item(45, 9).Value = 13079428  # I45
$ws.Cells.Item(45, 10).Value = 5000.778  # J45
$ws.Cells.Item(45, 11).Value = 13079428  # K45
$ws.Cells.Item(45, 12).Value = 5000.778  # L45
$ws.Cells.Item(45, 13).Value = -13079051  # M45
$ws.Cells.Item(45, 14).Value = -5754.778  # N45
$ws.Cells.Item(61, 8).Value = 3410.7896  # H61
$ws.Cells.Item(61, 9).Value = 3427.5  # I61
$ws.Cells.Item(61, 10).Value = 3321.6667  # J61
$ws.Cells.Item(61, 11).Value = 3427.5  # K61
$ws.Cells.Item(61, 12).Value = 3321.6667  # L61
$ws.Cells.Item(61, 13).Value = -3215.5  # M61
$ws.Cells.Item(61, 14).Value = -3745.6667  # N61
$ws.Cells.Item(74, 8).Value = 203178.8  # H74
$ws.Cells.Item(74, 9).Value = 87186.36  # I74
$ws.Cells.Item(74, 11).Value = 87186.36  # K74
$ws.Cells.Item(74, 13).Value = -86312.36  # M74
$ws.Cells.Item(77, 8).Value = 203178.8  # H77
$ws.Cells.Item(77, 9).Value = 87186.36  # I77
$ws.Cells.Item(77, 11).Value = 435931.8  # K77
$ws.Cells.Item(77, 13).Value = -431563.8  # M77
$ws.Cells.Item(98, 8).Value = 0  # H98
$ws.Cells.Item(98, 10).Value = 0  # J98
$ws.Cells.Item(98, 12).Value = 0  # L98
$ws.Cells.Item(98, 14).ClearContents()  # N98
$ws.Cells.Item(122, 8).Value = 633205.7  # H122
$ws.Cells.Item(122, 9).Value = 1826.9231  # I122
$ws.Cells.Item(122, 10).Value = 2978326.8  # J122
$ws.Cells.Item(122, 11).Value = 5480.7693  # K122
$ws.Cells.Item(122, 12).Value = 8934980.399999999  # L122
$ws.Cells.Item(122, 13).Value = -3030.7693  # M122
$ws.Cells.Item(122, 14).Value = -8939880.399999999  # N122
$ws.Cells.Item(132, 8).Value = 1918.7838  # H132
$ws.Cells.Item(132, 9).Value = 1411.7587  # I132
$ws.Cells.Item(132, 10).Value = 3756.75  # J132
$ws.Cells.Item(132, 11).Value = 4235.2761  # K132
$ws.Cells.Item(132, 12).Value = 11270.25  # L132
$ws.Cells.Item(132, 13).Value = -1705.2761  # M132
$ws.Cells.Item(132, 14).Value = -16330.25  # N132
$ws.Cells.Item(136, 8).Value = 3410.7896  # H136
$ws.Cells.Item(136, 9).Value = 3427.5  # I136
$ws.Cells.Item(136, 10).Value = 3321.6667  # J136
$ws.Cells.Item(136, 11).Value = 10282.5  # K136
$ws.Cells.Item(136, 12).Value = 9965.000100000001  # L136
$ws.Cells.Item(136, 13).Value = -7732.5  # M136
$ws.Cells.Item(136, 14).Value = -15065.0001  # N136

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2301.6316  # H20
$ws.Cells.Item(20, 9).Value = 1963.5385  # I20
$ws.Cells.Item(20, 10).Value = 3034.1667  # J20
$ws.Cells.Item(20, 11).Value = 1963.5385  # K20
$ws.Cells.Item(20, 12).Value = 3034.1667  # L20
$ws.Cells.Item(20, 13).Value = -1716.5385  # M20
$ws.Cells.Item(20, 14).Value = -3528.1667  # N20
$ws.Cells.Item(105, 8).Value = 3474904.8  # H105
$ws.Cells.Item(105, 9).Value = 4467049.5  # I105
$ws.Cells.Item(105, 10).Value = 2397.5  # J105
$ws.Cells.Item(105, 11).Value = 4467049.5  # K105
$ws.Cells.Item(105, 12).Value = 2397.5  # L105
$ws.Cells.Item(105, 13).Value = -4465302.5  # M105
$ws.Cells.Item(105, 14).Value = -5891.5  # N105

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 37783.5  # H31
$ws.Cells.Item(31, 9).Value = 1735.7273  # I31
$ws.Cells.Item(31, 10).Value = 64218.535  # J31
$ws.Cells.Item(31, 11).Value = 1735.7273  # K31
$ws.Cells.Item(31, 12).Value = 64218.535  # L31
$ws.Cells.Item(31, 13).Value = -1440.7273  # M31
$ws.Cells.Item(31, 14).Value = -64808.535  # N31
$ws.Cells.Item(34, 8).Value = 37783.5  # H34
$ws.Cells.Item(34, 9).Value = 1735.7273  # I34
$ws.Cells.Item(34, 10).Value = 64218.535  # J34
$ws.Cells.Item(34, 11).Value = 1735.7273  # K34
$ws.Cells.Item(34, 12).Value = 64218.535  # L34
$ws.Cells.Item(34, 13).Value = -1533.7273  # M34
$ws.Cells.Item(34, 14).Value = -64622.535  # N34
$ws.Cells.Item(58, 8).Value = 3665.4666  # H58
$ws.Cells.Item(58, 9).Value = 3663.6667  # I58
$ws.Cells.Item(58, 10).Value = 3668.1667  # J58
$ws.Cells.Item(58, 11).Value = 3663.6667  # K58
$ws.Cells.Item(58, 12).Value = 3668.1667  # L58
$ws.Cells.Item(58, 13).Value = -3460.6667  # M58
$ws.Cells.Item(58, 14).Value = -4074.1667  # N58
$ws.Cells.Item(99, 8).Value = 3473.4119  # H99
$ws.Cells.Item(99, 9).Value = 2618.7778  # I99
$ws.Cells.Item(99, 10).Value = 4434.875  # J99
$ws.Cells.Item(99, 11).Value = 2618.7778  # K99
$ws.Cells.Item(99, 12).Value = 4434.875  # L99
$ws.Cells.Item(99, 13).Value = -1120.7778  # M99
$ws.Cells.Item(99, 14).Value = -7430.875  # N99
$ws.Cells.Item(126, 8).Value = 3473.4119  # H126
$ws.Cells.Item(126, 9).Value = 2618.7778  # I126
$ws.Cells.Item(126, 10).Value = 4434.875  # J126
$ws.Cells.Item(126, 11).Value = 7856.3334  # K126
$ws.Cells.Item(126, 12).Value = 13304.625  # L126
$ws.Cells.Item(126, 13).Value = -5386.3334  # M126
$ws.Cells.Item(126, 14).Value = -18244.625  # N126
$ws.Cells.Item(132, 8).Value = 101508.37  # H132
$ws.Cells.Item(132, 9).Value = 64111  # I132
$ws.Cells.Item(132, 11).Value = 192333  # K132
$ws.Cells.Item(132, 13).Value = -189803  # M132
$ws.Cells.Item(136, 8).Value = 3665.4666  # H136
$ws.Cells.Item(136, 9).Value = 3663.6667  # I136
$ws.Cells.Item(136, 10).Value = 3668.1667  # J136
$ws.Cells.Item(136, 11).Value = 10991.0001  # K136
$ws.Cells.Item(136, 12).Value = 11004.5001  # L136
$ws.Cells.Item(136, 13).Value = -8441.000100000001  # M136
$ws.Cells.Item(136, 14).Value = -16104.5001  # N136

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 739.6  # H3
$ws.Cells.Item(3, 9).Value = 739.6  # I3
$ws.Cells.Item(3, 11).Value = 2218.8  # K3
$ws.Cells.Item(3, 13).Value = -2106.8  # M3
$ws.Cells.Item(56, 8).Value = 10006051  # H56
$ws.Cells.Item(56, 9).Value = 10006051  # I56
$ws.Cells.Item(56, 11).Value = 10006051  # K56
$ws.Cells.Item(56, 13).Value = -10005521  # M56
$ws.Cells.Item(133, 8).Value = 2026.2858  # H133
$ws.Cells.Item(133, 9).Value = 2026.2858  # I133
$ws.Cells.Item(133, 11).Value = 6078.857400000001  # K133
$ws.Cells.Item(133, 13).Value = -1018.857400000001  # M133
$ws.Cells.Item(134, 8).Value = 3965  # H134
$ws.Cells.Item(134, 9).Value = 3965  # I134
$ws.Cells.Item(134, 11).Value = 11895  # K134
$ws.Cells.Item(134, 13).Value = -6825  # M134

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 7998502.5  # H102
$ws.Cells.Item(102, 9).Value = 13890886  # I102
$ws.Cells.Item(102, 11).Value = 13890886  # K102
$ws.Cells.Item(102, 13).Value = -13889264  # M102
$ws.Cells.Item(113, 8).Value = 83334090  # H113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 14).ClearContents()  # N113
$ws.Cells.Item(122, 8).Value = 231042.31  # H122
$ws.Cells.Item(122, 9).Value = 298405.88  # I122
$ws.Cells.Item(122, 10).Value = 6497.1113  # J122
$ws.Cells.Item(122, 11).Value = 895217.64  # K122
$ws.Cells.Item(122, 12).Value = 19491.3339  # L122
$ws.Cells.Item(122, 13).Value = -892767.64  # M122
$ws.Cells.Item(122, 14).Value = -24391.3339  # N122
$ws.Cells.Item(126, 8).Value = 6135703.5  # H126
$ws.Cells.Item(126, 9).Value = 2843828  # I126
$ws.Cells.Item(126, 10).Value = 16669706  # J126
$ws.Cells.Item(126, 11).Value = 8531484  # K126
$ws.Cells.Item(126, 12).Value = 50009118  # L126
$ws.Cells.Item(126, 13).Value = -8529014  # M126
$ws.Cells.Item(126, 14).Value = -50014058  # N126
$ws.Cells.Item(132, 8).Value = 3188  # H132
$ws.Cells.Item(132, 9).Value = 2561.2222  # I132
$ws.Cells.Item(132, 11).Value = 7683.6666  # K132
$ws.Cells.Item(132, 13).Value = -5153.6666  # M132

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3668.2285  # H40
$ws.Cells.Item(40, 9).Value = 2888  # I40
$ws.Cells.Item(40, 11).Value = 2888  # K40
$ws.Cells.Item(40, 13).Value = -2752  # M40
$ws.Cells.Item(93, 8).Value = 27782926  # H93
$ws.Cells.Item(93, 9).Value = 33337710  # I93
$ws.Cells.Item(93, 10).Value = 8999.5  # J93
$ws.Cells.Item(93, 11).Value = 33337710  # K93
$ws.Cells.Item(93, 12).Value = 8999.5  # L93
$ws.Cells.Item(93, 13).Value = -33336462  # M93
$ws.Cells.Item(93, 14).Value = -11495.5  # N93
$ws.Cells.Item(100, 8).Value = 2930.4688  # H100
$ws.Cells.Item(100, 9).Value = 2877.8928  # I100
$ws.Cells.Item(100, 11).Value = 2877.8928  # K100
$ws.Cells.Item(100, 13).Value = -2336.8928  # M100
$ws.Cells.Item(122, 8).Value = 4358.5654  # H122
$ws.Cells.Item(122, 9).Value = 3114.389  # I122
$ws.Cells.Item(122, 10).Value = 8837.6  # J122
$ws.Cells.Item(122, 11).Value = 9343.167000000001  # K122
$ws.Cells.Item(122, 12).Value = 26512.8  # L122
$ws.Cells.Item(122, 13).Value = -6893.167000000001  # M122
$ws.Cells.Item(122, 14).Value = -31412.8  # N122
$ws.Cells.Item(132, 8).Value = 6601.851  # H132
$ws.Cells.Item(132, 9).Value = 7030.7095  # I132
$ws.Cells.Item(132, 10).Value = 5770.9375  # J132
$ws.Cells.Item(132, 11).Value = 21092.1285  # K132
$ws.Cells.Item(132, 12).Value = 17312.8125  # L132
$ws.Cells.Item(132, 13).Value = -18562.1285  # M132
$ws.Cells.Item(132, 14).Value = -22372.8125  # N132

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 23810580  # H81
$ws.Cells.Item(81, 9).Value = 83333830  # I81
$ws.Cells.Item(81, 10).Value = 1279  # J81
$ws.Cells.Item(81, 11).Value = 166667660  # K81
$ws.Cells.Item(81, 12).Value = 2558  # L81
$ws.Cells.Item(81, 13).Value = -166666599  # M81
$ws.Cells.Item(81, 14).Value = -4680  # N81
$ws.Cells.Item(84, 8).Value = 23810580  # H84
$ws.Cells.Item(84, 9).Value = 83333830  # I84
$ws.Cells.Item(84, 10).Value = 1279  # J84
$ws.Cells.Item(84, 11).Value = 833338300  # K84
$ws.Cells.Item(84, 12).Value = 12790  # L84
$ws.Cells.Item(84, 13).Value = -833332996  # M84
$ws.Cells.Item(84, 14).Value = -23398  # N84
$ws.Cells.Item(100, 8).Value = 994.5  # H100
$ws.Cells.Item(100, 9).Value = 1229.2  # I100
$ws.Cells.Item(100, 10).Value = 603.3333  # J100
$ws.Cells.Item(100, 11).Value = 2458.4  # K100
$ws.Cells.Item(100, 12).Value = 1206.6666  # L100
$ws.Cells.Item(100, 13).Value = -1917.4  # M100
$ws.Cells.Item(100, 14).Value = -2288.6666  # N100
$ws.Cells.Item(107, 8).Value = 52635176  # H107
$ws.Cells.Item(107, 9).Value = 58827356  # I107
$ws.Cells.Item(107, 11).Value = 176482068  # K107
$ws.Cells.Item(107, 13).Value = -176480148  # M107
$ws.Cells.Item(126, 8).Value = 2302.2666  # H126
$ws.Cells.Item(126, 9).Value = 2519.5  # I126
$ws.Cells.Item(126, 10).Value = 1433.3334  # J126
$ws.Cells.Item(126, 11).Value = 7558.5  # K126
$ws.Cells.Item(126, 12).Value = 4300.0002  # L126
$ws.Cells.Item(126, 13).Value = -5088.5  # M126
$ws.Cells.Item(126, 14).Value = -9240.0002  # N126
